# P-W compatible with tap changers
#
# Bus 1 (sheet "Busos") becomes a tap-changer / PV bus: new P/Q figures, a
# voltage setpoint, and its "P" entry keeps the quote-prefixed number style
# that the source workbook carries on that cell. Bus 2 gets revised P/Q and
# loses its shunt entry, and the old bus-3 row is cleared out (its data was
# folded into bus 1). The "Topologia" (line) sheet gets revised R/X/B and
# tap-ratio figures, condensing from 4 data rows to 3.

$wb = $excel.ActiveWorkbook

$busos = $wb.Worksheets.Item("Busos")
$topologia = $wb.Worksheets.Item("Topologia")

# ---------------------------------------------------------------------------
# Sheet "Busos"
# ---------------------------------------------------------------------------

# Slack bus voltage setpoint tweak
$busos.Range("D2").Value = 1.1000000000000001

# Bus 1 becomes a PV (tap-changer) bus: new P, Q and voltage setpoint.
$busos.Range("B3").Value = -2
$busos.Range("C3").Value = -0.5
$busos.Range("D3").Value = 1
$busos.Range("F3").Value = "PV"

# B3's "P" value is quote-prefixed in the source file (number stored with a
# leading apostrophe). Build that style on a scratch cell and copy just the
# formatting across, then wipe the scratch cell again.
$busos.Range("E3").Value = "'-2"
$busos.Range("E3").Copy() | Out-Null
$busos.Range("B3").PasteSpecial(-4122) | Out-Null
$busos.Range("E3").Clear() | Out-Null

# Bus 2 keeps its PQ type but gets new P/Q figures; its shunt cell is cleared.
$busos.Range("B4").Value = -1
$busos.Range("C4").Value = -0.3
$busos.Range("G4").ClearContents() | Out-Null

# Remove old bus 3 row (its data now folded into bus 1 above)
$busos.Range("A5:G5").ClearContents() | Out-Null

$busos.Range("G5").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Topologia"
# ---------------------------------------------------------------------------

$topologia.Range("C2").Value = 0
$topologia.Range("D2").Value = 0.04
$topologia.Range("E2").Value = 0.1
$topologia.Range("F2").Value = 0.95

$topologia.Range("A3").Value = 1
$topologia.Range("C3").Value = 0.01
$topologia.Range("D3").Value = 0.070000000000000007
$topologia.Range("E3").Value = 0.04

$topologia.Range("A4").Value = 0
$topologia.Range("B4").Value = 2
$topologia.Range("C4").Value = 0.04
$topologia.Range("D4").Value = 0.09
$topologia.Range("E4").Value = 0
$topologia.Range("F4").Value = 1

# Remove old row 5 (topology now condensed to 3 data rows)
$topologia.Range("A5:F5").ClearContents() | Out-Null

$topologia.Activate() | Out-Null
$topologia.Range("F5").Select() | Out-Null
